$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 9.54645530522332
$ws.Range("G2").Value = 9.170167643570377
$ws.Range("H2").Value = 9.933142240451385
$ws.Range("I2").Value = 0.002217226290430173
$ws.Range("J2").Value = 0.001960015122806
$ws.Range("K2").Value = 0.002514946819293703
$ws.Range("L2").Value = 0.008321064069526204
$ws.Range("M2").Value = 0.008099617639449332
$ws.Range("N2").Value = 0.008554356793520618

$ws.Range("F3").Value = 0.04688919838484605
$ws.Range("G3").Value = 0.0466187969992323
$ws.Range("H3").Value = 0.04717019893536182
$ws.Range("I3").Value = 0.04532124684966276
$ws.Range("J3").Value = 0.04506091196824184
$ws.Range("K3").Value = 0.04559082519707216
$ws.Range("L3").Value = 0.04694087289506213
$ws.Range("M3").Value = 0.04667050585486655
$ws.Range("N3").Value = 0.04722190739345494

$ws.Range("F4").Value = 9.593344503608165
$ws.Range("G4").Value = 9.216786440569606
$ws.Range("H4").Value = 9.980312439386747
$ws.Range("I4").Value = 0.04753847314009294
$ws.Range("J4").Value = 0.04702092709104784
$ws.Range("K4").Value = 0.04810577201636586
$ws.Range("L4").Value = 0.05526193696458832
$ws.Range("M4").Value = 0.05477012349431588
$ws.Range("N4").Value = 0.05577626418697556
